# CU-1vhwepw added ntc for api Register
# Add two new "negative test case" rows to the apiTest sheet, each
# re-using the existing Srdjan/Rados identity but with a tweaked
# username/password, plus their own hyperlink on the email cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("apiTest")

# --- Row 3: email with a malformed hyperlink display text ("...htecgroup") ---
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:srdjan.rados@htecgroup", "", "", "srdjan.rados@htecgroup")
$ws.Range("A3").Value = "srdjan.rados@htecgroup.com"
$ws.Range("A3").Font.Underline = 0
$ws.Range("A3").Font.ColorIndex = 1
$ws.Range("B3").Value = "Qwertysha1@"
$ws.Range("C3").Value = "Srdjan"
$ws.Range("D3").Value = "Rados"
$ws.Rows("3").RowHeight = 13.8

# --- Row 4: email with correct hyperlink, different (invalid) password ---
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:srdjan.rados@htecgroup.com", "", "", "srdjan.rados@htecgroup.com")
$ws.Range("A4").Value = "srdjan.rados@htecgroup.com"
$ws.Range("A4").Font.Underline = 0
$ws.Range("A4").Font.ColorIndex = 1
$ws.Range("B4").Value = "Qwertysha"
$ws.Range("C4").Value = "Srdjan"
$ws.Range("D4").Value = "Rados"
$ws.Rows("4").RowHeight = 13.8

# Restore the selection to match the author's final cursor position.
$ws.Range("C20").Select()
